$d = $word.ActiveDocument
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 0) Remove the "_GoBack" bookmark from its old location (end of the "Nowe
#    klasy specjalistyczne..." paragraph) first, before any other structural
#    edit, to avoid any ambiguity from having two same-named bookmarks.
# ---------------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
[void]$old.Delete()

# ---------------------------------------------------------------------------
# 1) Paragraph "utrata polaczenia z BD w trakcie dzialania aplikacji - ..."
#    -> strike-through the whole paragraph (incl. the paragraph mark), drop
#    the " - initialize(true)?" tail, and re-split the remaining text into
#    two runs.
# ---------------------------------------------------------------------------
$target = $d.Paragraphs.Item(2)
$full = $target.Range
$full.Font.StrikeThrough = 1

$full2 = $d.Paragraphs.Item(2).Range
$contentRange = $d.Range($full2.Start, $full2.End - 1)

$runsXml = $pkgOpen + '<w:p>' + `
  '<w:r><w:rPr><w:strike/><w:color w:val="00B0F0"/></w:rPr><w:t>utrata po' + [char]0x0142 + [char]0x0105 + 'czenia z B</w:t></w:r>' + `
  '<w:r><w:rPr><w:strike/><w:color w:val="00B0F0"/></w:rPr><w:t>D w trakcie dzia' + [char]0x0142 + 'ania aplikacji</w:t></w:r>' + `
  '</w:p>' + $pkgClose
[void]$contentRange.InsertXML($runsXml)

# ---------------------------------------------------------------------------
# 2) The following empty paragraph becomes the new home of the "_GoBack"
#    bookmark.
# ---------------------------------------------------------------------------
$blank = $d.Paragraphs.Item(3).Range
$bmXml = $pkgOpen + '<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + $pkgClose
[void]$blank.InsertXML($bmXml)

# ---------------------------------------------------------------------------
# 3) Insert the new "Statyczna klasa do obslugi MessageBox'ow" paragraph
#    right after "Nowe klasy specjalistyczne...".
# ---------------------------------------------------------------------------
$afterRange = $d.Paragraphs.Item(16).Range
[void]$afterRange.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(17).Range
$newXml = $pkgOpen + '<w:p>' + `
  '<w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/><w:color w:val="00B0F0"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:strike/><w:color w:val="00B0F0"/></w:rPr><w:t xml:space="preserve">Statyczna klasa do obs' + [char]0x0142 + 'ugi </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:strike/><w:color w:val="00B0F0"/></w:rPr><w:t>MessageBox' + [char]0x2019 + [char]0x00f3 + 'w</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>' + $pkgClose
[void]$newPara.InsertXML($newXml)

Write-Host "edit complete"
